# Refresh the cryptocurrency symbol list (price + 1h volume columns) with the
# latest scraped values. Source cells are stored as literal text (inlineStr)
# so we force a Text number format before writing, otherwise Excel would
# auto-coerce numeric-looking strings (e.g. "261.24") or percentages
# (e.g. "1.08%") into numbers and rewrite their formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "261.24"
    "E2" = "1.08%"
    "D3" = "27.10"
    "E3" = "0.80%"
    "D4" = "4.700"
    "E4" = "1.21%"
    "D5" = "0.06190"
    "E5" = "3.52%"
    "D6" = "6.683"
    "E6" = "0.60%"
    "D7" = "0.8503"
    "E7" = "-0.65%"
    "D8" = "0.9160"
    "E8" = "-0.43%"
    "E9" = "1.37%"
    "E10" = "0.40%"
    "D11" = "0.07092"
    "E11" = "1.08%"
    "D12" = "0.03153"
    "E12" = "3.29%"
    "D13" = "0.09043"
    "E13" = "-0.76%"
    "D14" = "0.001528"
    "E14" = "0.11%"
    "D15" = "0.0006175"
    "E15" = "2.33%"
    "D16" = "0.006117"
    "E16" = "-1.29%"
    "E17" = "0.26%"
    "D18" = "3.176"
    "E18" = "0.80%"
    "E19" = "-1.21%"
    "E21" = "0.87%"
    "D22" = "4.077"
    "E22" = "0.97%"
    "D23" = "0.04217"
    "E23" = "-0.40%"
    "D24" = "0.001216"
    "E24" = "-0.10%"
    "E25" = "-5.52%"
    "E26" = "0.02%"
    "D27" = "0.0001577"
    "E27" = "-7.84%"
    "D40" = "0.03892"
    "E40" = "1.68%"
    "E41" = "-0.15%"
    "E42" = "8.88%"
    "E43" = "8.49%"
    "D44" = "0.002183"
    "E44" = "-10.13%"
    "D45" = "0.00005156"
    "E45" = "1.09%"
    "E46" = "0.04%"
    "D48" = "0.1667"
    "E48" = "43.06%"
    "D49" = "0.00002100"
    "E49" = "0.04%"
    "D50" = "0.0002000"
    "E50" = "0.04%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
